$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells before writing so that
# numeric-looking values (e.g. "252.57") are stored as literal
# text, matching the inlineStr cells in the source workbook.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '99.370.35'
$ws.Range('E2').Value = '  +0.93%  '
$ws.Range('D3').Value = '3.282.39'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '252.57'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '621.00'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('E7').Value = '  +20.13%  '
$ws.Range('D8').Value = '0.398'
$ws.Range('E8').Value = '  +3.38%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '0.962'
$ws.Range('E10').Value = '  +19.99%  '
$ws.Range('D11').Value = '3.278.78'
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('D12').Value = '0.199'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Value = '39.16'
$ws.Range('E13').Value = '  +9.29%  '
$ws.Range('D14').Value = '99.100.71'
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').Value = '0.0000246'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '3.879.04'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').Value = '5.45'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '3.276.17'
$ws.Range('E18').Value = '  -2.39%  '
$ws.Range('D19').Value = '3.42'
$ws.Range('E19').Value = '  -4.42%  '
$ws.Range('D20').Value = '15.22'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').Value = '6.31'
$ws.Range('E21').Value = '  +7.85%  '
$ws.Range('D22').Value = '486.16'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '9.26'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('D24').Value = '0.0000199'
$ws.Range('E24').Value = '  -3.10%  '
$ws.Range('D25').Value = '5.61'
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('D26').Value = '88.91'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('D27').Value = '0.318'
$ws.Range('E27').Value = '  +28.48%  '
$ws.Range('D28').Value = '11.93'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').Value = '3.431.59'
$ws.Range('E29').Value = '  -2.83%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').Value = '0.136'
$ws.Range('E31').Value = '  +8.51%  '
$ws.Range('E32').Value = '  +2.51%  '
$ws.Range('D33').Value = '10.31'
$ws.Range('E33').Value = '  +11.28%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('D35').Value = '27.77'
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('D36').Value = '0.470'
$ws.Range('E36').Value = '  +5.27%  '
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').Value = '7.14'
$ws.Range('E38').Value = '  -3.14%  '
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('D40').Value = '24.82'
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').Value = '485.54'
$ws.Range('E41').Value = '  -5.73%  '
$ws.Range('D42').Value = '3.62'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Value = '1.22'
$ws.Range('E43').Value = '  -3.00%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '0.767'
$ws.Range('E45').Value = '  -1.17%  '
$ws.Range('D46').Value = '3.06'
$ws.Range('E46').Value = '  -5.81%  '
$ws.Range('D47').Value = '1.93'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').Value = '157.17'
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('D49').Value = '0.845'
$ws.Range('E49').Value = '  +7.08%  '
$ws.Range('D50').Value = '7.23'
$ws.Range('E50').Value = '  +14.60%  '
$ws.Range('D51').Value = '4.69'
$ws.Range('E51').Value = '  +4.09%  '

# Reset the cell style back to Normal so no stray number-format
# style survives on the Price cells (keeps styles.xml clean).
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
